# Apply the edits described by the commit:
#  1. Rename header columns (A1:D1) to short English field names.
#  2. Title-case the Spanish connector words (de, del, el, la, las, los, y)
#     wherever they appear as standalone words inside the "state"/"municipality"
#     text columns (A and B).
#  3. Remove the trailing footnote rows (1214-1218) and shrink the used range
#     back down to A1:D1212.

function Fix-Text($s) {
    $s = $s -replace '\bde\b','De'
    $s = $s -replace '\bdel\b','Del'
    $s = $s -replace '\bel\b','El'
    $s = $s -replace '\blas\b','Las'
    $s = $s -replace '\bla\b','La'
    $s = $s -replace '\blos\b','Los'
    $s = $s -replace '\by\b','Y'
    return $s
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row before the footnote block that will be dropped.
$lastDataRow = 1212
# Rows holding the trailing footnote / metadata text that must be removed.
$firstFootnoteRow = 1214
$lastFootnoteRow = 1218

# 1) Re-title-case the municipality/state name columns (A and B), rows 2..lastDataRow.
for ($r = 2; $r -le $lastDataRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null -and $a -ne "") {
        $ws.Cells.Item($r, 1).Value = (Fix-Text $a)
    }
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null -and $b -ne "") {
        $ws.Cells.Item($r, 2).Value = (Fix-Text $b)
    }
}

# 2) Replace the header row with short machine-friendly field names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Delete the trailing footnote rows so the sheet (and its dimension) ends at row 1212.
$footnoteRange = $ws.Range("A" + $firstFootnoteRow + ":A" + $lastFootnoteRow)
$footnoteRange.EntireRow.Delete()
